$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowRange($rangeAddr, $vals) {
    $n = $vals.Length
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0,$i] = $vals[$i]
    }
    $ws.Range($rangeAddr).Value = $arr
}

Set-RowRange "B2:F2" @(1.02,1.007517189882853,1.010163023887673,1.010173097929085,1.005706564766229)
Set-RowRange "I2:N2" @(1.023594999628091,1.012788370573424,1.013034046953649,1.013044090525428,1.008591131989706,1.008322484023764)
Set-RowRange "B3:F3" @(1.02,1.008777202253876,1.011333579330225,1.011251613067968,1.007621938967665)
Set-RowRange "I3:N3" @(1.023504579208683,1.013678231869978,1.014007589187392,1.013925852672849,1.010306393046803,1.008628179638942)
Set-RowRange "B4:F4" @(1.02,1.009591924433953,1.012090770948062,1.011949264209806,1.008860251563744)
Set-RowRange "I4:N4" @(1.023443335729026,1.014253038152566,1.014636742348009,1.014495612020765,1.011414844666994,1.008825338359081)
Set-RowRange "B5:F5" @(1.02,1.009934297215639,1.012409041674543,1.012242507248352,1.009380597362823)
Set-RowRange "I5:N5" @(1.023416932628374,1.014494452321081,1.014901051889395,1.014734950459729,1.011880505766299,1.008908069990505)
Set-RowRange "B6:F6" @(1.02,1.00999177525264,1.012462477755981,1.012291741235623,1.009467951983473)
Set-RowRange "I6:N6" @(1.023412460900047,1.014534973177195,1.014945419794748,1.014775125478151,1.011958673196528,1.00892195199838)
Set-RowRange "B7:F7" @(1.02,1.009596499765874,1.012095023899101,1.011953182728163,1.008867205379284)
Set-RowRange "I7:N7" @(1.023442985510381,1.014256264857257,1.014640274793233,1.014498810809175,1.01142106814292,1.008826444426729)
Set-RowRange "B8:F8" @(1.02,1.007943142064123,1.010558669789743,1.010537634729582,1.006354102649742)
Set-RowRange "I8:N8" @(1.023565007130014,1.013089311152642,1.013363226625048,1.01334225360938,1.009171117304427,1.008425929781417)
Set-RowRange "B9:F9" @(1.02,1.005024977690948,1.007849435501134,1.008041425703835,1.001916947197381)
Set-RowRange "I9:N9" @(1.023759139049229,1.011025225612804,1.011106664175804,1.011297989151956,1.00519486974853,1.007715175401641)
Set-RowRange "B10:F10" @(1.02,1.003076043972148,1.006041700691044,1.006375848360975,0.9989521099072133)
Set-RowRange "I10:N10" @(1.023874603872166,1.009643729766757,1.009597866282022,1.009930740491634,1.002535517254903,1.007237918939261)
Set-RowRange "B11:F11" @(1.02,1.002231237537959,1.005258501728547,1.005654247396487,0.9976665156659683)
Set-RowRange "I11:N11" @(1.023921309042974,1.009044187569629,1.008943439922908,1.009337617796226,1.001381798509034,1.007030436511823)
Set-RowRange "B12:F12" @(1.02,1.001917296819953,1.004967516128487,1.005386149240515,0.9971887018993026)
Set-RowRange "I12:N12" @(1.023938164268129,1.008821284461941,1.008700186052415,1.009117137245418,1.000952910667623,1.006953242781246)
Set-RowRange "B13:F13" @(1.02,1.001984644691728,1.005029936802072,1.005443660141311,0.997291207817982)
Set-RowRange "I13:N13" @(1.023934571058996,1.008869107346712,1.008752372633972,1.009164438753571,1.001044924427577,1.006969806798484)
Set-RowRange "B14:F14" @(1.02,1.002205290031676,1.005234450237169,1.005632087636969,0.9976270253773073)
Set-RowRange "I14:N14" @(1.023922712353276,1.009025766565623,1.00892333599205,1.009319396262326,1.001346353635605,1.007024058221722)
Set-RowRange "B15:F15" @(1.02,1.002341217912532,1.005360448208953,1.005748175508192,0.9978338951597572)
Set-RowRange "I15:N15" @(1.023915340507297,1.009122261995733,1.00902864939447,1.009414848270602,1.001532027921628,1.007057467644408)
Set-RowRange "B16:F16" @(1.02,1.003132091409819,1.006093669388492,1.00642372997058,0.9990373912417067)
Set-RowRange "I16:N16" @(1.023871434968932,1.00968349069489,1.009641274701102,1.009970080692238,1.002612038143036,1.007251671330573)
Set-RowRange "B17:F17" @(1.02,1.003627938879378,1.006553479023196,1.006847378814025,0.9997918190738039)
Set-RowRange "I17:N17" @(1.023843013862359,1.010035171399794,1.010025258865606,1.010318067441776,1.003288900350847,1.007373267829813)
Set-RowRange "B18:F18" @(1.02,1.003917071380771,1.00682163628692,1.007094448185242,1.000231691490648)
Set-RowRange "I18:N18" @(1.023826118601784,1.010240171289094,1.010249123836221,1.010520936836497,1.003683491195755,1.007444113279416)
Set-RowRange "B19:F19" @(1.02,1.004015643525776,1.006913064020034,1.007178686161864,1.000381647839342)
Set-RowRange "I19:N19" @(1.023820303822167,1.010310049115234,1.010325438077179,1.010590092297764,1.003818001053591,1.007468256260234)
Set-RowRange "B20:F20" @(1.02,1.003574748174463,1.006504150171866,1.006801929245292,0.9997108940757602)
Set-RowRange "I20:N20" @(1.023846096030143,1.009997452816174,1.009984072046662,1.010280742685765,1.003216301416195,1.007360229932181)
Set-RowRange "B21:F21" @(1.02,1.002140319446499,1.005174228116869,1.005576602219859,0.9975281435260565)
Set-RowRange "I21:N21" @(1.02392621804881,1.008979640063875,1.008872996294841,1.009273769817298,1.001257599897156,1.007008086004692)
Set-RowRange "B22:F22" @(1.02,1.001237612647247,1.004337642242208,1.004805821458388,0.9961540939351058)
Set-RowRange "I22:N22" @(1.023973741382756,1.00833850422362,1.008173428531071,1.008639669918192,1.000024081236138,1.006785952305349)
Set-RowRange "B23:F23" @(1.02,1.001716234765744,1.004781172838824,1.005214463276155,0.9968826672004006)
Set-RowRange "I23:N23" @(1.023948818288664,1.008678497332355,1.008544378016596,1.00897591209132,1.000678187788766,1.006903778894003)
Set-RowRange "B24:F24" @(1.02,1.003598783033373,1.006526439886278,1.006822466066262,0.9997474611321396)
Set-RowRange "I24:N24" @(1.023844704313072,1.010014496620955,1.010002682926594,1.010297608462272,1.003249106403019,1.007366121443385)
Set-RowRange "B25:F25" @(1.02,1.005779985857745,1.008550097602331,1.008686993834474,1.00306518222948)
Set-RowRange "I25:N25" @(1.02371141742819,1.011559783695783,1.01169080246928,1.011827242669469,1.006224271044895,1.007899520777064)

Write-Output "applied vm_pu updates for 24 rows"